$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# New week 17 betting lines to append below the existing data (last existing row is 225)
$data = @(
    @(17, "DAL_WAS", 47.5, 5.5),
    @(17, "DET_MIN", 46.5, -1.5),
    @(17, "DEN_KC", 44.5, 4.5),
    @(17, "HOU_LAC", 43.5, 3.5),
    @(17, "BAL_GB", 46.5, -1.5),
    @(17, "ARI_CIN", 47.5, 3.5),
    @(17, "SEA_CAR", 43.5, -1.5),
    @(17, "PIT_CLE", 40.5, -2.5),
    @(17, "NO_TEN", 42.5, 3.5),
    @(17, "TB_MIA", 47.5, -1.5),
    @(17, "JAX_IND", 45.5, 1.5),
    @(17, "NE_NYJ", 41.5, -1.5),
    @(17, "NYG_LV", 42.5, 2.5),
    @(17, "PHI_BUF", 46.5, 1.5),
    @(17, "CHI_SF", 46.5, 3.5),
    @(17, "LA_ATL", 46.5, -3)
)

$startRow = 226
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Update the view to reflect the newly scrolled/selected position
$ws.Range("I229").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 220

Write-Host "Added $($data.Count) rows starting at row $startRow"
